$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.648.97'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.748.75'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.03%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.11'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.44'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.749.34'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.23%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.95%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.29'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.19'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000247'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.374.85'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.94%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.749.20'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.693.46'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.49%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.33%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.08'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.47%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.92'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +20.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '494.72'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.47%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000153'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +10.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.31'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.08%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.87%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.62%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.23'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.44%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.47%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +7.73%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.90'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.89'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.92%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.894.08'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.97%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.683.45'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.78%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.25%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.02'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.82%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.84'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.52%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.64%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '440.58'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.84'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.61%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.86'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.28%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.84%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.21'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.61'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.802.37'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.75%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.38%  '
